$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "245.30"
Set-TextValue "E2" "-0.16%"
Set-TextValue "D3" "26.53"
Set-TextValue "E3" "4.24%"
Set-TextValue "D4" "5.120"
Set-TextValue "E4" "-0.18%"
Set-TextValue "D5" "0.05584"
Set-TextValue "E5" "-0.09%"
Set-TextValue "D6" "6.479"
Set-TextValue "E6" "-0.78%"
Set-TextValue "D7" "0.8172"
Set-TextValue "E7" "-0.01%"
Set-TextValue "D8" "0.8393"
Set-TextValue "E8" "-0.08%"
Set-TextValue "D9" "0.1331"
Set-TextValue "E9" "-0.82%"
Set-TextValue "D10" "0.06991"
Set-TextValue "E10" "0.63%"
Set-TextValue "D11" "0.02871"
Set-TextValue "E11" "1.08%"
Set-TextValue "D12" "0.09382"
Set-TextValue "E12" "-0.03%"
Set-TextValue "D13" "0.001515"
Set-TextValue "E13" "-0.22%"
Set-TextValue "D14" "0.0005971"
Set-TextValue "E14" "0.32%"
Set-TextValue "D15" "0.006141"
Set-TextValue "E15" "0.02%"
Set-TextValue "D16" "3.633"
Set-TextValue "E16" "3.65%"
Set-TextValue "D17" "3.039"
Set-TextValue "E17" "0.72%"
Set-TextValue "D20" "0.03067"
Set-TextValue "E20" "-2.95%"
Set-TextValue "D21" "0.1299"
Set-TextValue "E21" "-2.21%"
Set-TextValue "D22" "3.739"
Set-TextValue "E22" "-0.18%"
Set-TextValue "D23" "0.04584"
Set-TextValue "E23" "-2.99%"
Set-TextValue "E24" "2.45%"
Set-TextValue "E25" "-0.13%"
Set-TextValue "D26" "0.004515"
Set-TextValue "E26" "5.79%"
Set-TextValue "D27" "0.00009598"
Set-TextValue "E27" "-1.10%"
Set-TextValue "D28" "0.0001396"
Set-TextValue "E28" "0.56%"
Set-TextValue "E40" "-0.61%"
Set-TextValue "D41" "0.1370"
Set-TextValue "E41" "1.65%"
Set-TextValue "D42" "0.002620"
Set-TextValue "E42" "-0.22%"
Set-TextValue "E43" "-44.44%"
Set-TextValue "D44" "0.008213"
Set-TextValue "E44" "-1.11%"
Set-TextValue "D45" "0.00005337"
Set-TextValue "E45" "0.73%"
Set-TextValue "E46" "-0.02%"
Set-TextValue "E47" "-51.56%"
Set-TextValue "D48" "0.002552"
Set-TextValue "E48" "20.39%"
Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "-0.02%"
Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "-0.02%"
